$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update id_rombel (column B) values for rows 2-25: alternate 6/8 by row parity
for ($r = 2; $r -le 25; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Cells.Item($r, 2).Value = 6
    } else {
        $ws.Cells.Item($r, 2).Value = 8
    }
}

# Update the frozen pane / view position and selection
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
